$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.307.65'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.671.87'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.83%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.34'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.54'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.670.48'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.82%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.525'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("E10").Value = '  -5.46%  '
$ws.Range("E11").Value = '  -4.59%  '
$ws.Range("E12").Value = '  -4.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.42'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.84%  '
$ws.Range("E14").Value = '  -5.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.282.32'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.670.83'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.335.30'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.61'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +5.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.15'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.18%  '
$ws.Range("E20").Value = '  -3.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '492.93'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("E22").Value = '  -4.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.718'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.36'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("E25").Value = '  -6.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000137'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.14'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.07%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.96'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.86%  '
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("E31").Value = '  -6.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.65'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.44'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.806.28'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.04%  '
$ws.Range("E35").Value = '  -5.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.606.82'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.87%  '
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.986'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.69%  '
$ws.Range("E39").Value = '  -5.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.130'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.322'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '432.59'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -10.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.62'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("E44").Value = '  -6.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.76'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.89%  '
$ws.Range("E46").Value = '  -2.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.39'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -7.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '142.72'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.750.87'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0346'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.05%  '
